$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("DSD")

# --- Indicator row (row 5): codelist changed from CL_INDICATOR to the new
#     fishing-method-indicators codelist, and it is no longer a "CDCL" item.
$ws.Range("F5").Value = "CL_FISHING_METHOD_INDICATORS"
$ws.Range("G5").Value = "N"

# --- Fishing-method dimension rows (8-12): each used to reference its own
#     per-method codelist (CL_FISHING_GLEANING / _LINE / _NET / _SPEAR /
#     _OTHER_METHOD); now they all share one Yes/No codelist, and the CDCL
#     flag flips on.
$ws.Range("F8").Value  = "CL_COM_YESNO"
$ws.Range("G8").Value  = "Y"

$ws.Range("F9").Value  = "CL_COM_YESNO"
$ws.Range("G9").Value  = "Y"

$ws.Range("F10").Value = "CL_COM_YESNO"
$ws.Range("G10").Value = "Y"

$ws.Range("F11").Value = "CL_COM_YESNO"
$ws.Range("G11").Value = "Y"

$ws.Range("F12").Value = "CL_COM_YESNO"
$ws.Range("G12").Value = "Y"

# --- Column F widened to fit the longer codelist name.
$ws.Columns("F:F").AutoFit()

# --- DSD becomes the active sheet / tab, with F12 selected.
$ws.Activate() | Out-Null
$ws.Range("F12").Select() | Out-Null
